$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.690.19'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '2.551.69'
$ws.Range('E3').Value = '  -2.59%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.95'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.42'
$ws.Range('E6').Value = '  -1.37%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.520'
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('D9').Value = '2.551.52'
$ws.Range('E9').Value = '  -2.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.165'
$ws.Range('E10').Value = '  -2.38%  '
$ws.Range('E11').Value = '  -1.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.355'
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('D14').Value = '3.011.78'
$ws.Range('E14').Value = '  -2.82%  '
$ws.Range('D15').Value = '70.535.34'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000178'
$ws.Range('E16').Value = '  -5.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.46'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('D18').Value = '2.546.01'
$ws.Range('E18').Value = '  -2.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.93'
$ws.Range('E19').Value = '  +1.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.45'
$ws.Range('E20').Value = '  -5.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '355.57'
$ws.Range('E21').Value = '  -4.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.95'
$ws.Range('E22').Value = '  -2.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.05'
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.21'
$ws.Range('E25').Value = '  -1.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.06'
$ws.Range('E26').Value = '  -2.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.19'
$ws.Range('E27').Value = '  -1.23%  '
$ws.Range('D28').Value = '2.657.87'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '0.0₃0925'
$ws.Range('E30').Value = '  -1.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.97'
$ws.Range('E31').Value = '  +0.58%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.29'
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '474.43'
$ws.Range('E33').Value = '  -2.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.77'
$ws.Range('E34').Value = '  -0.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('E36').Value = '  +3.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.71'
$ws.Range('E37').Value = '  -1.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.04'
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.70'
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.88'
$ws.Range('E41').Value = '  +0.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.324'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.62'
$ws.Range('E43').Value = '  -4.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.37'
$ws.Range('E44').Value = '  -6.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.17'
$ws.Range('E45').Value = '  -13.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.42'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '145.52'
$ws.Range('E47').Value = '  -2.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.539'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.55'
$ws.Range('E49').Value = '  -2.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.63'
$ws.Range('E50').Value = '  -1.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0739'
$ws.Range('E51').Value = '  -0.26%  '
